$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.05619466666666667
$ws.Range("I2").Value = 0.04986276087265156
$ws.Range("J2").Value = 0.07297477932340853
$ws.Range("M2").Value = 10.250695
$ws.Range("N2").Value = 30.752085
$ws.Range("O2").Value = 0.2501330709220828
$ws.Range("P2").Value = 0.259830908271274
$ws.Range("Q2").Value = 0.5760343886266667
$ws.Range("R2").Value = 5.18430949764
$ws.Range("S2").Value = 0.01247232550172981
$ws.Range("T2").Value = 0.01896110319249702

# Row 3
$ws.Range("G3").Value = 0.05619466666666667
$ws.Range("I3").Value = 0.04986276087265156
$ws.Range("J3").Value = 0.07297477932340853
$ws.Range("O3").Value = 0.2460808482884365
$ws.Range("P3").Value = 0.2556215780794031
$ws.Range("Q3").Value = 0.5667024775013334
$ws.Range("S3").Value = 0.01227027049354556
$ws.Range("T3").Value = 0.01865392825064588

# Row 4
$ws.Range("G4").Value = 0.05619466666666667
$ws.Range("I4").Value = 0.04986276087265156
$ws.Range("J4").Value = 0.07297477932340853
$ws.Range("M4").Value = 7.311799000000001
$ws.Range("N4").Value = 21.935397
$ws.Range("O4").Value = 0.1784193889131434
$ws.Range("P4").Value = 0.1853368357235283
$ws.Range("Q4").Value = 0.4108841075386667
$ws.Range("R4").Value = 3.697956967848
$ws.Range("S4").Value = 0.008896483324420687
$ws.Range("T4").Value = 0.0135249146874233

# Row 5
$ws.Range("G5").Value = 0.05619466666666667
$ws.Range("I5").Value = 0.04986276087265156
$ws.Range("J5").Value = 0.07297477932340853
$ws.Range("M5").Value = 4.588677499999999
$ws.Range("N5").Value = 9.177354999999999
$ws.Range("O5").Value = 0.1119709438770801
$ws.Range("P5").Value = 0.07754142475796089
$ws.Range("Q5").Value = 0.2578592025533333
$ws.Range("R5").Value = 1.54715521532
$ws.Range("S5").Value = 0.005583180399227933
$ws.Range("T5").Value = 0.005658568360134882

# Row 6
$ws.Range("G6").Value = 0.05619466666666667
$ws.Range("I6").Value = 0.04986276087265156
$ws.Range("J6").Value = 0.07297477932340853
$ws.Range("M6").Value = 8.745164000000001
$ws.Range("N6").Value = 26.235492
$ws.Range("O6").Value = 0.2133957479992572
$ws.Range("P6").Value = 0.2216692531678338
$ws.Range("Q6").Value = 0.4914315759253334
$ws.Range("R6").Value = 4.422884183328001
$ws.Range("S6").Value = 0.01064050115372757
$ws.Range("T6").Value = 0.01617626483270745

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.070792
$ws.Range("H7").Value = 2.141584
$ws.Range("I7").Value = 0.9501372391273485
$ws.Range("J7").Value = 0.9270252206765914
$ws.Range("M7").Value = 10.250695
$ws.Range("N7").Value = 30.752085
$ws.Range("O7").Value = 0.2501330709220828
$ws.Range("P7").Value = 0.259830908271274
$ws.Range("Q7").Value = 10.97636220044
$ws.Range("R7").Value = 65.85817320263999
$ws.Range("S7").Value = 0.237660745420353
$ws.Range("T7").Value = 0.2408698050787769

# Row 8
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.070792
$ws.Range("H8").Value = 2.141584
$ws.Range("I8").Value = 0.9501372391273485
$ws.Range("J8").Value = 0.9270252206765914
$ws.Range("O8").Value = 0.2460808482884365
$ws.Range("P8").Value = 0.2556215780794031
$ws.Range("Q8").Value = 10.798542197752
$ws.Range("R8").Value = 64.79125318651199
$ws.Range("S8").Value = 0.233810577794891
$ws.Range("T8").Value = 0.2369676498287572

# Row 9
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.070792
$ws.Range("H9").Value = 2.141584
$ws.Range("I9").Value = 0.9501372391273485
$ws.Range("J9").Value = 0.9270252206765914
$ws.Range("M9").Value = 7.311799000000001
$ws.Range("N9").Value = 21.935397
$ws.Range("O9").Value = 0.1784193889131434
$ws.Range("P9").Value = 0.1853368357235283
$ws.Range("Q9").Value = 7.829415874808
$ws.Range("R9").Value = 46.976495248848
$ws.Range("S9").Value = 0.1695229055887227
$ws.Range("T9").Value = 0.171811921036105

# Row 10
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.070792
$ws.Range("H10").Value = 2.141584
$ws.Range("I10").Value = 0.9501372391273485
$ws.Range("J10").Value = 0.9270252206765914
$ws.Range("M10").Value = 4.588677499999999
$ws.Range("N10").Value = 9.177354999999999
$ws.Range("O10").Value = 0.1119709438770801
$ws.Range("P10").Value = 0.07754142475796089
$ws.Range("Q10").Value = 4.913519157579999
$ws.Range("R10").Value = 19.65407663032
$ws.Range("S10").Value = 0.1063877634778522
$ws.Range("T10").Value = 0.071882856397826

# Row 11
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.070792
$ws.Range("H11").Value = 2.141584
$ws.Range("I11").Value = 0.9501372391273485
$ws.Range("J11").Value = 0.9270252206765914
$ws.Range("M11").Value = 8.745164000000001
$ws.Range("N11").Value = 26.235492
$ws.Range("O11").Value = 0.2133957479992572
$ws.Range("P11").Value = 0.2216692531678338
$ws.Range("Q11").Value = 9.364251649888001
$ws.Range("R11").Value = 56.185509899328
$ws.Range("S11").Value = 0.2027552468455296
$ws.Range("T11").Value = 0.2054929883351263
